$d = $word.ActiveDocument

# Change 1: "a check box for remember me" -> "a check box for the option remember me"
$d.Content.Find.Execute(
    "a check box for remember me,link for forgotten password",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "a check box for the option remember me,link for forgotten password",
    2) | Out-Null

# Change 2: merge the three runs describing the SAVE button / VEHICLE_CATEGORY text
# into a single contiguous run with the same text (collapsing the run split).
$d.Content.Find.Execute(
    "and a SAVE button.The user will add new VEHICLE_CATEGORY for each driving cathegory he has in his driver licence.  The app will only allow  one driver/user. After the SAVE button is clicked the ADD_NEW_DRIVER button will be greyed out",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "and a SAVE button.The user will add new VEHICLE_CATEGORY for each driving cathegory he has in his driver licence.  The app will only allow  one driver/user. After the SAVE button is clicked the ADD_NEW_DRIVER button will be greyed out",
    2) | Out-Null
